# This workbook contains weekly price records for "Betarraga" (beet) at
# "Vega Monumental Concepción". A new weekly record (two rows: Primera and
# Segunda quality) was added right before the existing row that used to be
# row 297, pushing all subsequent rows down by two positions and extending
# the used range from A1:R422 to A1:R424.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 297-298; Excel will shift rows 297:422 down to
# 299:424 and automatically grow the sheet dimension to A1:R424.
$ws.Rows("297:298").Insert()

# Populate the newly inserted row 297 ("Primera" quality).
$ws.Range("A297").Value = 11
$ws.Range("B297").Value = "Vega Monumental Concepción"
$ws.Range("C297").Value = "Bíobío"
$ws.Range("D297").Value = 44924
$ws.Range("E297").Value = 8
$ws.Range("F297").Value = 100114014
$ws.Range("G297").Value = "Betarraga"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 600
$ws.Range("K297").Value = 700
$ws.Range("L297").Value = 800
$ws.Range("M297").Value = 750
$ws.Range("N297").Value = "$/paquete 5 unidades"
$ws.Range("O297").Value = "Región Metropolitana"
$ws.Range("P297").Value = 150
$ws.Range("Q297").Value = 5
$ws.Range("R297").Value = "Hortaliza"

# Populate the newly inserted row 298 ("Segunda" quality).
$ws.Range("A298").Value = 11
$ws.Range("B298").Value = "Vega Monumental Concepción"
$ws.Range("C298").Value = "Bíobío"
$ws.Range("D298").Value = 44924
$ws.Range("E298").Value = 8
$ws.Range("F298").Value = 100114014
$ws.Range("G298").Value = "Betarraga"
$ws.Range("H298").Value = "Sin especificar"
$ws.Range("I298").Value = "Segunda"
$ws.Range("J298").Value = 300
$ws.Range("K298").Value = 600
$ws.Range("L298").Value = 600
$ws.Range("M298").Value = 600
$ws.Range("N298").Value = "$/paquete 5 unidades"
$ws.Range("O298").Value = "Región Metropolitana"
$ws.Range("P298").Value = 120
$ws.Range("Q298").Value = 5
$ws.Range("R298").Value = "Hortaliza"
